$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Structural changes: drop the R1-2 resistor row (it gets merged into the
#    RN1 row below, becoming "RN1-2"), then re-insert a blank row where the
#    old row-13 separator used to be so a new "Total" summary row can live
#    at row 12 (pushing the separator + solder rows back down to 13-16).
# ---------------------------------------------------------------------------
$ws.Rows(7).Delete()
$ws.Rows(12).Insert()

# Drop the "Qty/panel" column (old G) entirely - board qty is used directly
# in the cost formula now instead of a panelised quantity.
$ws.Columns("G").Delete()

# ---------------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "Digikey PN"
$ws.Range("H1").Value = "Digikey Cost"
$ws.Range("I1").Value = "Cost/Board"

# ---------------------------------------------------------------------------
# 3. Per-row content fixes
# ---------------------------------------------------------------------------

# Row 3: C3-4 -> C3
$ws.Range("A3").Value = "C3"

# Row 7: RN1 -> RN1-2, qty 1 -> 2 (absorbs the deleted R1-2 resistor)
$ws.Range("A7").Value = "RN1-2"
$ws.Range("F7").Value = 2

# Row 8: "S1, JP3" -> "S1"
$ws.Range("A8").Value = "S1"

# Row 9: U1 voltage regulator part swap
$ws.Range("B9").Value = "MIC39100-5.0WS"
$ws.Range("D9").Value = "SOT-223 "
$ws.Range("G9").Value = "576-1173-ND"
$ws.Range("H9").Value = 1.66

# ---------------------------------------------------------------------------
# 4. Cost-per-board formulas (now board qty * unit cost, was panel qty *
#    unit cost through the dropped column)
# ---------------------------------------------------------------------------
$ws.Range("I2").Formula = "=`$F2*`$H2"
$ws.Range("I3:I11").FormulaR1C1 = "=RC6*RC8"

# ---------------------------------------------------------------------------
# 5. New Total row (row 12)
# ---------------------------------------------------------------------------
$ws.Range("G12").Value = $null
$ws.Range("H12").Value = "Total:"
$ws.Range("I12").Formula = "=SUM(I2:I11)"

# Bold currency-style label, right aligned, no border (new xf #6 in styles.xml)
$ws.Range("H12").Style = "Currency"
$ws.Range("H12").Font.Bold = $true
$ws.Range("H12").Borders.LineStyle = -4142      # xlLineStyleNone - no border
$ws.Range("H12").HorizontalAlignment = -4152    # xlRight

$ws.Range("I13").Value = $null
$ws.Range("I14").Value = $null
$ws.Range("I15").Value = $null
$ws.Range("I16").Value = $null

# ---------------------------------------------------------------------------
# 6. Selection cursor, matching the author's final click
# ---------------------------------------------------------------------------
$ws.Range("D20").Select()
